# string distance metric was added
# A new scraped listing was inserted as row 4, pushing the former rows 4-6 down
# to rows 5-7. Ranking (column A), refreshed tracking URLs and "time left"
# strings were updated across the sheet as a result of the re-scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4 so the former rows 4,5,6 become 5,6,7.
$ws.Rows.Item(4).Insert()

# Copy the formatting used by column A (rank) cells into the new row's A cell,
# matching the formatting already applied to the other data rows of column A.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)

# ---- Row 2 (unchanged item, refreshed rank + tracking URL) ----
$ws.Cells.Item(2, 1).Value = 3
$ws.Cells.Item(2, 3).Value = "https://www.ebay.com/itm/225299351425?hash=item3474e34781:g:WagAAOSwxeNjmfpZ&amdata=enc%3AAQAHAAAAoMZmeAwO9uGAOepXocLGZPSvZ7HVASv%2Bz9HfSuZcidVwwN8gERo1QG8THm0o8LhT0Z2CWTLpiihlW2Jqo9UGKd4chTi6G%2BhgFVtqmYv8EhdiAwzzeQ1HWHx7vheckA0u7oQz1RMfzVFLVKCVy1LGoFt4KCrHWtXYMEREawLxKhbfEQrJMXQKdzq1zPLiBr1GaiD7oXE3O67fN8jUf4ZFJSA%3D%7Ctkp%3ABk9SR7K8v8WjYQ"

# ---- Row 3 (unchanged item, refreshed rank + tracking URL) ----
$ws.Cells.Item(3, 1).Value = 5
$ws.Cells.Item(3, 3).Value = "https://www.ebay.com/itm/404050551629?epid=8047515616&hash=item5e134a234d:g:cYsAAOSwWPljk3xw&amdata=enc%3AAQAHAAAAoGiTEUPW5VtRzKGDLjO84d5VQscWIRHpASDiRaUr3QjG%2FuifRtUKslCBpUxzH3JynE4T6%2FM94mhy0pgeZkU2Wm30%2B5tOXwuKJtEfVusYws3B6WMTIpJGRtcUr8n8sMhRFD47g%2FySNakl3%2BhpFyuV6dq8fZdfnhWGTlUTQcF2xkE2BW2dxVkK8hYZJ5TRdiEXW18yq4SvssZ9jfiIgPSdN0o%3D%7Ctkp%3ABk9SR7K8v8WjYQ"

# ---- Row 4 (brand-new listing inserted by the re-scrape) ----
$ws.Cells.Item(4, 1).Value = 0
$ws.Cells.Item(4, 2).Value = "Novo anúncioPowerColor Red Devil Radeon RX 6900 XT 16 Go GDDR6 Carte Graphique TBE"
$ws.Cells.Item(4, 3).Value = "https://www.ebay.com/itm/394381873231?epid=7050529243&hash=item5bd2fdd04f:g:vzsAAOSwyghjnKug&amdata=enc%3AAQAHAAAAkB3NEuFpGgXCau7OnPW1AXiUui0xbUF3H3rVBhVWvtl69O%2B%2FsoxuZY0TjiRlsciqGAvEUTjBveVTwFlRKIaOvLwMMvG2%2FBb1tHrIxUyipe8Tyf6kIXbJ1WjS%2FVAtyFxa2ex7s2Z%2F30CsUXjtNr4Z7FUMQq2dLzOxBACG1BQ0b00vpfK5TCflICHJuBUaDXT%2BPQ%3D%3D%7Ctkp%3ABk9SR7K8v8WjYQ"
$ws.Cells.Item(4, 4).Value = "Seminovo"
$ws.Cells.Item(4, 5).Value = 580
$ws.Cells.Item(4, 6).Value = "de França"
$ws.Cells.Item(4, 7).Value = 29
$ws.Cells.Item(4, 8).Value = 609
$ws.Cells.Item(4, 9).Value = "Yes"
$ws.Cells.Item(4, 10).Value = "No other purchase options"
$ws.Cells.Item(4, 11).Value = "2d 22h restantes"

# ---- Row 5 (was row 4 before insertion; rank + URL + price/location/shipping/time refreshed) ----
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 3).Value = "https://www.ebay.com/itm/295403696780?epid=7050209760&hash=item44c76eb28c:g:bw4AAOSwChJjk1Kw&amdata=enc%3AAQAHAAAAkBAKU7HgZCqFBzuU0e3plvMFS2not8HtJB9dLimFOyJn5Ux9KmyJYM0gh6SCd8i7oryvdHANMHqRmXDvbeQPqGVnsI%2F0mIoayRInFFnyzvkmv8V%2BZtlTWEHSmeZltJIS4A9kQiaa6tYAN8VPc8CXfXOlmDe7eQtSJ9FNzuc3SO4jJyrYqrUqwuDvNYnCaGRwTw%3D%3D%7Ctkp%3ABk9SR7K8v8WjYQ"
$ws.Cells.Item(5, 5).Value = 685
$ws.Cells.Item(5, 6).Value = "de Alemanha"
$ws.Cells.Item(5, 7).Value = 15
$ws.Cells.Item(5, 8).Value = 700
$ws.Cells.Item(5, 11).Value = "2d 20h restantes"

# ---- Row 6 (was row 5 before insertion; rank + URL + time refreshed) ----
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 3).Value = "https://www.ebay.com/itm/295418574962?hash=item44c851b872:g:QQkAAOSwJNhjmuak&amdata=enc%3AAQAHAAAAoInV0dSL5sn815Yx6uFdOd3lkWn1srcTHg1wPp96pvDmNNaWXTEEAxNLjwP9W9oQXcDT0RyBToeg%2FXguFDjMKh61kysiv8FudGZ4TPL7849J8YN8i1ebZSHO%2BXMoNE2ri%2FkWe%2FYGl%2F8%2BJarWYhv%2BOmXZjdh6GtFk0HEhlV51XQLfC5WpxY1gUCmYtnZo5V04bdg0zXRWcWlCZz0%2FUs61xG8%3D%7Ctkp%3ABk9SR7K8v8WjYQ"
$ws.Cells.Item(6, 11).Value = "5d 14h restantes"

# ---- Row 7 (was row 6 before insertion; rank + refreshed tracking URL + time) ----
$ws.Cells.Item(7, 1).Value = 2
$ws.Cells.Item(7, 3).Value = "https://www.ebay.com/itm/285074976656?epid=14050520857&hash=item425fcaef90:g:EGcAAOSwdJBjlNZH&amdata=enc%3AAQAHAAAAoJj3syuSee2Bl6pg78DOlpmggJb6JULGmdUQt%2FVrBsLHPL9HUcAVRC8O%2Bv3DNiZK4gh2osDg6uCCMXJHIKf%2BdCwMMcY4hXD%2B033hIViPqboTGecvmPYv%2FZ4DxejHXvK%2F9C43EbKmkxW22VB1%2BI3UUduEhzy9LM3JkgJudCPFYuuiYDlvMGfTG%2FP7yy3ifhhQGmY2rOhiZfv%2Bj2sFKIisJNA%3D%7Ctkp%3ABk9SR7K8v8WjYQ"
$ws.Cells.Item(7, 11).Value = "1d 23h restantes"
